$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A95").Value = "2024-11-05 00:00:00"
$ws.Range("B95").Value = 74950
$ws.Range("C95").Value = 10520.18
$ws.Range("D95").Value = 9309.9
$ws.Range("E95").Value = 7.1039
